# UI-AddCategory.data.xlsx — "Modified Add category and TrnasferNote script"
#
# 1) Insert a new row at 78 (parentCategory.enter.values / the option[3] xpath)
#    right after the existing "select.ParentCategory" row, pushing everything
#    below down by one.
# 2) Append 5 new label/value rows at the end of the data (after the old last
#    row, "spec.add.delete.button") for the new specification / delete-popup
#    related locators.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) insert new row under "ParentCategory" select control ---------------
$ws.Rows("78:78").Insert()

# copy formatting (style ids 1/2 used by column A/B) down from the row that
# was just pushed to 79, so the new row matches the sheet's normal look
$ws.Range("A79:B79").Copy()
$ws.Range("A78:B78").PasteSpecial(-4122)

$ws.Range("A78").Value = "parentCategory.enter.values"
$ws.Range("B78").Value = "//main[@class='mb-5']/section/div/div[2]/form/fieldset[2]/div/div/select/option[3]"

# --- 2) append new rows for specification popup / delete locators ----------
$lastRow = 148  # original last row of data ("spec.add.delete.button"), now at 149

$newRows = @(
    @("popup.message.xpath", "//div[@class='b-toaster-slot vue-portal-target']/div/div/div"),
    @("specification.title.values.xpath", "//main[@class='mb-5']/section/div/div[2]/form/table[1]/tbody/tr/td[2]"),
    @("specification.button.xpath", "//main[@class='mb-5']/section/div/div[2]/form/table[1]/tbody/tr/td[3]"),
    @("delete.specification.button", "/div/button[2]"),
    @("Selection.field", "/div/button[2]//div2")
)

$row = $lastRow + 2  # account for the row inserted in step 1 (148 -> 149), next free row is 150
foreach ($pair in $newRows) {
    $ws.Range("A79:B79").Copy()
    $ws.Range("A" + $row + ":B" + $row).PasteSpecial(-4122)
    $ws.Range("A" + $row).Value = $pair[0]
    $ws.Range("B" + $row).Value = $pair[1]
    $row = $row + 1
}

# land the selection on the next empty row, same as Excel would after
# manually keying in the last new row of data
$ws.Range("A" + $row).Select() | Out-Null

Write-Output "done"
